# Auto-generated edits applying updated market price data to Sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 67666.60000000001
$ws.Range("I2").Value = 466.33334
$ws.Range("J2").Value = 112466.78
$ws.Range("K2").Value = 466.33334
$ws.Range("L2").Value = 112466.78
$ws.Range("M2").Value = -353.33334
$ws.Range("N2").Value = -112692.78

$ws.Range("H9").Value = 8699.833000000001
$ws.Range("I9").Value = 304.5
$ws.Range("K9").Value = 304.5
$ws.Range("M9").Value = -135.5

$ws.Range("H15").Value = 736.6316
$ws.Range("I15").Value = 736.6316
$ws.Range("K15").Value = 2209.8948
$ws.Range("M15").Value = -2040.8948

$ws.Range("H40").Value = 3963.5557
$ws.Range("J40").Value = 6974.25
$ws.Range("L40").Value = 6974.25
$ws.Range("N40").Value = -7324.25

$ws.Range("H43").Value = 2799
$ws.Range("I43").Value = 2323.5
$ws.Range("K43").Value = 2323.5
$ws.Range("M43").Value = -2254.5

$ws.Range("H53").Value = 3789.6667
$ws.Range("J53").Value = 682
$ws.Range("L53").Value = 682
$ws.Range("N53").Value = -1956

$ws.Range("H58").Value = 751.6
$ws.Range("I58").Value = 751.6
$ws.Range("K58").Value = 2254.8
$ws.Range("M58").Value = -2104.8

$ws.Range("H86").Value = 3176.111
$ws.Range("J86").Value = 2933.3333
$ws.Range("L86").Value = 2933.3333
$ws.Range("N86").Value = -5179.3333

$ws.Range("H89").Value = 3176.111
$ws.Range("J89").Value = 2933.3333
$ws.Range("L89").Value = 14666.6665
$ws.Range("N89").Value = -25898.6665

$ws.Range("H111").Value = 2808
$ws.Range("J111").Value = 2754
$ws.Range("L111").Value = 8262
$ws.Range("N111").Value = -14396

$ws.Range("H131").Value = 8661.666999999999
$ws.Range("J131").Value = 14999
$ws.Range("L131").Value = 44997
$ws.Range("N131").Value = -55077

$ws.Range("H138").Value = 9560536
$ws.Range("I138").Value = 4809379.5
$ws.Range("J138").Value = 25001794
$ws.Range("K138").Value = 14428138.5
$ws.Range("L138").Value = 75005382
$ws.Range("M138").Value = -14422998.5
$ws.Range("N138").Value = -75015662

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3137.9836
$ws.Range("I32").Value = 2845.1086
$ws.Range("K32").Value = 2845.1086
$ws.Range("M32").Value = -2558.1086

$ws.Range("H61").Value = 3376
$ws.Range("I61").Value = 2548.1516
$ws.Range("J61").Value = 4617.773
$ws.Range("K61").Value = 2548.1516
$ws.Range("L61").Value = 4617.773
$ws.Range("M61").Value = -2336.1516
$ws.Range("N61").Value = -5041.773

$ws.Range("H136").Value = 3376
$ws.Range("I136").Value = 2548.1516
$ws.Range("J136").Value = 4617.773
$ws.Range("K136").Value = 7644.4548
$ws.Range("L136").Value = 13853.319
$ws.Range("M136").Value = -5094.4548
$ws.Range("N136").Value = -18953.319

$ws.Range("H138").Value = 174678.62
$ws.Range("J138").Value = 174678.62
$ws.Range("L138").Value = 174678.62
$ws.Range("N138").Value = -184958.62

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1888.2188
$ws.Range("I31").Value = 1163.5
$ws.Range("J31").Value = 2612.9375
$ws.Range("K31").Value = 1163.5
$ws.Range("L31").Value = 2612.9375
$ws.Range("M31").Value = -868.5
$ws.Range("N31").Value = -3202.9375

$ws.Range("H34").Value = 1888.2188
$ws.Range("I34").Value = 1163.5
$ws.Range("J34").Value = 2612.9375
$ws.Range("K34").Value = 1163.5
$ws.Range("L34").Value = 2612.9375
$ws.Range("M34").Value = -961.5
$ws.Range("N34").Value = -3016.9375

$ws.Range("H115").Value = 47229.168
$ws.Range("J115").Value = 47229.168
$ws.Range("L115").Value = 47229.168
$ws.Range("N115").Value = -49579.168

$ws.Range("H134").Value = 5410.875
$ws.Range("I134").Value = 4381.222
$ws.Range("K134").Value = 13143.666
$ws.Range("M134").Value = -10608.666

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws.Range("H138").Value = 80000
$ws.Range("J138").Value = 80000
$ws.Range("L138").Value = 80000
$ws.Range("N138").Value = -90280

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 1008
$ws.Range("J26").Value = 80.333336
$ws.Range("L26").Value = 241.000008
$ws.Range("N26").Value = -817.000008

$ws.Range("H55").Value = 76936640
$ws.Range("J55").Value = 90924820
$ws.Range("L55").Value = 272774460
$ws.Range("N55").Value = -272774814

$ws.Range("H92").Value = 202.27272
$ws.Range("I92").Value = 114.4
$ws.Range("K92").Value = 343.2
$ws.Range("M92").Value = 904.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5184.1934
$ws.Range("I102").Value = 4689
$ws.Range("K102").Value = 4689
$ws.Range("M102").Value = -3067

$ws.Range("H126").Value = 166670200
$ws.Range("I126").Value = 166670200
$ws.Range("K126").Value = 500010600
$ws.Range("M126").Value = -500008130

$ws.Range("H132").Value = 1863.225
$ws.Range("I132").Value = 1634.8064
$ws.Range("J132").Value = 2650
$ws.Range("K132").Value = 4904.4192
$ws.Range("L132").Value = 7950
$ws.Range("M132").Value = -2374.4192
$ws.Range("N132").Value = -13010

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3588.375
$ws.Range("I22").Value = 2016.4286
$ws.Range("K22").Value = 2016.4286
$ws.Range("M22").Value = -1721.4286

$ws.Range("H27").Value = 3588.375
$ws.Range("I27").Value = 2016.4286
$ws.Range("K27").Value = 2016.4286
$ws.Range("M27").Value = -1909.4286

$ws.Range("H40").Value = 28409.354
$ws.Range("I40").Value = 29684.938
$ws.Range("K40").Value = 29684.938
$ws.Range("M40").Value = -29548.938

$ws.Range("H61").Value = 17789
$ws.Range("I61").Value = 15678.417
$ws.Range("J61").Value = 30452.5
$ws.Range("K61").Value = 15678.417
$ws.Range("L61").Value = 30452.5
$ws.Range("M61").Value = -15476.417
$ws.Range("N61").Value = -30856.5

$ws.Range("H82").Value = 1423.9524
$ws.Range("J82").Value = 1385.4286
$ws.Range("L82").Value = 1385.4286
$ws.Range("N82").Value = -2107.4286

$ws.Range("H85").Value = 1423.9524
$ws.Range("J85").Value = 1385.4286
$ws.Range("L85").Value = 1385.4286
$ws.Range("N85").Value = -3881.4286

$ws.Range("H113").Value = 17789
$ws.Range("I113").Value = 15678.417
$ws.Range("J113").Value = 30452.5
$ws.Range("K113").Value = 15678.417
$ws.Range("L113").Value = 30452.5
$ws.Range("M113").Value = -13508.417
$ws.Range("N113").Value = -34792.5

$ws.Range("I136").Value = 2664.4092
$ws.Range("J136").Value = 66669800
$ws.Range("K136").Value = 7993.2276
$ws.Range("L136").Value = 200009400
$ws.Range("M136").Value = -5443.2276
$ws.Range("N136").Value = -200014500

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3599.5
$ws.Range("I122").Value = 3429.8462
$ws.Range("K122").Value = 10289.5386
$ws.Range("M122").Value = -7839.5386

$ws.Range("H126").Value = 66670180
$ws.Range("I126").Value = 83336720
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 250010160
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -250007690
$ws.Range("N126").Value = -16940
